# Fix Training Data Issue (#48)
# The "Date" column (BF) held the source filename-derived string
# "6-21-2013-14" for every data row. That's off by a day versus how the
# NBA stats site actually labelled the date, so correct it to the real
# ISO date "2014-06-21" for every row in the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "6-21-2013-14"
$newDate = "2014-06-21"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $firstCol + $used.Columns.Count - 1

# Locate the "Date" header cell in row 1 so this keeps working even if
# the column shifts around.
$dateCol = 0
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Value2
    if ($header -eq "Date") {
        $dateCol = $c
        break
    }
}

if ($dateCol -eq 0) {
    # Fallback to the known column (BF) from the original layout.
    $dateCol = 58
}

# Force the target cells to Text so Excel doesn't reinterpret the
# replacement string "2014-06-21" as a date serial number.
$dataRange = $ws.Range($ws.Cells.Item($firstRow + 1, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
$dataRange.NumberFormat = "@"

for ($row = $firstRow + 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldDate) {
        $cell.Value = $newDate
    }
}
